# Issue #5: stock data output to json file
# Adds a "property_category" column to the 股票 (stock) sheet and fixes a
# handful of shared-string typos (stray internal spaces / full-width commas).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)   # 股票 (stock) sheet

# --- Fix stray spaces / full-width punctuation inside existing text values ---
$ws.Range("B2").Value = "★新光合成纖維股份有限公司"
$ws.Range("B3").Value = "★鴻海精密工業股份有限公司"
$ws.Range("B8").Value = "★健喬信元醫藥生技股份有限公司"
$ws.Range("G8").Value = "1746410"

# --- Insert a new "property_category" column between "total" and "date" ---
$ws.Range("H1").EntireColumn.Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H8").Value = "stock"
